# [Fonds de solidarite] Add 2022-06-10 data
# Update "nombre_aides" (column C) and "montant_total" (column E) for the
# rows whose underlying source data changed with the new 2022-06-10 extract.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 45;  C = 23378;  E = 99234004 },
    @{ Row = 47;  C = 3598;   E = 31499398 },
    @{ Row = 54;  C = 60294;  E = 353615499 },
    @{ Row = 61;  C = 51959;  E = 133917851 },
    @{ Row = 92;  C = 409173; E = 1595662657 },
    @{ Row = 94;  C = 94212;  E = 918364469 },
    @{ Row = 95;  C = 50780;  E = 933227656 },
    @{ Row = 97;  C = 2161;   E = 214282109 },
    @{ Row = 104; C = 135248; E = 272240066 },
    @{ Row = 141; C = 80475;  E = 280728664 },
    @{ Row = 142; C = 168976; E = 681799746 },
    @{ Row = 154; C = 201571; E = 786779728 },
    @{ Row = 182; C = 71;     E = 11214004 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
